# daily auto push: 2025-10-06 13:37 UTC
# Append the new data-log row (row 71) to the bottom of the table on Sheet1:
#   A71 = "2025/10/06"  (text - leading apostrophe stops Excel's autodetect
#                         from turning the slash-separated string into a date
#                         serial, matching every other row in column A)
#   B71 = "月"           (text)
#   C71 = 21             (number)
#   D71 = 201            (number)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A71").Value = "'2025/10/06"
$ws.Range("B71").Value = "月"
$ws.Range("C71").Value = 21
$ws.Range("D71").Value = 201
